$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Grupp" values in column B for rows 7-19 (6a was not correct)
$ws.Range("B7").Value = 9
$ws.Range("B8").Value = 8
$ws.Range("B9").Value = 7
$ws.Range("B10").Value = 6
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 4
$ws.Range("B13").Value = 3
$ws.Range("B14").Value = 2
$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = 9
$ws.Range("B18").Value = 8
$ws.Range("B19").Value = 7

# Update the active selection to match the author's final cursor position
$ws.Range("B20").Select()
